# Remove the "Winter Facts" slide (slide 2).
# The presentation originally had 3 slides:
#   1. Christmas Presentation (title)
#   2. Winter Facts (bullet list)
#   3. Christmas Tree (picture)
# After the edit, the "Winter Facts" slide is deleted entirely, so the
# "Christmas Tree" slide becomes the new slide 2.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$s.Delete()
